$d = $word.ActiveDocument
$p = $d.Paragraphs.Item(4)
$r = $p.Range
$r.Find.Execute("40 to 48", $false, $false, $false, $false, $false, $true, 1, $false, "35 to 45", 2)
Write-Host "step1:" $d.Paragraphs.Item(4).Range.Text
$r2 = $d.Paragraphs.Item(4).Range
$r2.Find.Execute("witch is usually between 8 and 10 rounds", $false, $false, $false, $false, $false, $true, 1, $false, "which is 8 or 12. You may try other numbers but make sure there are no surrogate matches (a 1 in the column with all of the zeros)", 2)
Write-Host "step2:" $d.Paragraphs.Item(4).Range.Text
